$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.171.12'
$ws.Range("E2").Value = '  -4.81%  '
$ws.Range("D3").Value = '2.879.01'
$ws.Range("E3").Value = '  -5.44%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '484.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -6.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.87'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.19%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.414'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -6.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.06'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -5.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.102'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.343'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -6.98%  '
$ws.Range("D12").Value = '3.379.18'
$ws.Range("E12").Value = '  -5.69%  '
$ws.Range("E13").Value = '  -4.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.37'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000155'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -8.53%  '
$ws.Range("D16").Value = '55.227.50'
$ws.Range("E16").Value = '  -4.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.90'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.69%  '
$ws.Range("D18").Value = '2.882.57'
$ws.Range("E18").Value = '  -5.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.27'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -6.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '309.25'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -7.97%  '
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.77'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.473'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '61.65'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("E27").Value = '  -6.33%  '
$ws.Range("D28").Value = '0.0₃0827'
$ws.Range("E28").Value = '  -13.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.26'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -9.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.85'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -8.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.71'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.44'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.68%  '
$ws.Range("E33").Value = '  -9.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.33'
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.31'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -9.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.50'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -7.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '24.12'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.16'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -9.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0642'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -7.01%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.84'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.72%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.624'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.58'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -8.05%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.080.14'
$ws.Range("E44").Value = '  -10.35%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.31'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -9.56%  '
$ws.Range("B46").Value = 'Cosmos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.76'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.88%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.896'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -10.90%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0226'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.81%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.30'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.95%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0831'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -7.56%  '
$ws.Range("B51").Value = 'ZEEBU'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.93'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.37%  '
